$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update Status column (E) to "Concluido" for the tasks that were finished.
$rowsToComplete = @(4, 6, 7, 10, 11, 14, 15, 18, 19, 21, 22, 23, 26, 27, 29, 30)
foreach ($r in $rowsToComplete) {
    $ws.Cells.Item($r, 5).Value = "Concluido"
}

# Move the active selection to F22, matching the saved cursor position.
$ws.Range("F22").Select()
